$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-08 20:27:22", 0.0016),
    @("2023-12-08 20:28:52", 0.0058),
    @("2023-12-08 20:29:25", 0.0022)
)

$startRow = 130
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
